# "added buttons for prediction"
# Relabel the Customer Segment categories and record rows that were
# (re)classified by the new prediction buttons, then update the sheet's
# view/selection to reflect where the user last clicked.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename two segment labels everywhere they occur in the AE (Customer
#    Segment) column - this keeps every row that used to read
#    "Highclass Families" / "Upper Middle Class" in sync with the rename.
for ($r = 2; $r -le 66; $r++) {
    $cell = $ws.Cells.Item($r, 31)
    $val = $cell.Value()
    if ($val -eq "Highclass Families") {
        $cell.Value = "Family"
    } elseif ($val -eq "Upper Middle Class") {
        $cell.Value = "Upper Class"
    }
}

# 2) Rows that the new prediction buttons reclassified out of "No Segment"
#    into a concrete customer segment.
$predictions = @{
    23 = "Middle Class"
    30 = "Middle Class"
    31 = "Middle Class"
    32 = "Middle Class"
    38 = "Family"
    41 = "Middle Class"
    43 = "Upper Class"
    44 = "Upper Class"
    45 = "Upper Class"
    48 = "Highclass Youngsters"
    50 = "Middle Class"
    51 = "Middle Class"
    53 = "Middle Class"
    54 = "Middle Class"
    57 = "Middle Class"
    60 = "Middle Class"
    61 = "Middle Class"
    63 = "Middle Class"
    64 = "Middle Class"
    65 = "Middle Class"
    66 = "Upper Class"
}

foreach ($row in $predictions.Keys) {
    $ws.Cells.Item($row, 31).Value = $predictions[$row]
}

# 3) Update the view to show where the work happened - scrolled to U15,
#    with AF32 as the active cell/selection.
$excel.Goto($ws.Range("U15"), $true)
$win = $excel.ActiveWindow
$win.ScrollRow = 15
$win.ScrollColumn = 21
$ws.Range("AF32").Select()
